$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain number (e.g. "7.10", "0.629") need to be
# forced to Text format first, otherwise Excel auto-converts them to numeric values
# and silently drops formatting such as trailing zeros, losing the original text look.
$ws.Range('D2').Value = '42.183.11'
$ws.Range('E2').Value = '  +2.19%  '
$ws.Range('D3').Value = '2.223.03'
$ws.Range('E3').Value = '  +1.48%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '252.04'
$ws.Range('E5').Value = '  -1.43%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.624'
$ws.Range('E6').Value = '  -0.11%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '68.22'
$ws.Range('E7').Value = '  -0.50%  '
$ws.Range('E8').Value = '  -0.06%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.629'
$ws.Range('E9').Value = '  +7.75%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '39.17'
$ws.Range('E10').Value = '  +2.78%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '59.86'
$ws.Range('E11').Value = '  +2.82%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0942'
$ws.Range('E12').Value = '  +0.03%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '7.10'
$ws.Range('E13').Value = '  -0.66%  '
$ws.Range('E14').Value = '  +0.15%  '
$ws.Range('D15').Value = '2.556.99'
$ws.Range('E15').Value = '  +1.54%  '
$ws.Range('B16').Value = 'Polygon'
$ws.Range('C16').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.879'
$ws.Range('E16').Value = '  +0.94%  '
$ws.Range('B17').Value = 'Chainlink'
$ws.Range('C17').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '14.65'
$ws.Range('E17').Value = '  +0.51%  '
$ws.Range('D18').Value = '2.219.36'
$ws.Range('E18').Value = '  +1.59%  '
$ws.Range('D19').Value = '42.090.96'
$ws.Range('E19').Value = '  +2.15%  '
$ws.Range('D20').Value = '0.0₃0965'
$ws.Range('E20').Value = '  +1.75%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.20'
$ws.Range('E21').Value = '  -0.17%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '72.70'
$ws.Range('E22').Value = '  +1.33%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '232.39'
$ws.Range('E23').Value = '  -0.27%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.04'
$ws.Range('E24').Value = '  -2.23%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.89'
$ws.Range('E25').Value = '  +0.17%  '
$ws.Range('E26').Value = '  -4.18%  '
$ws.Range('E27').Value = '  +0.28%  '
$ws.Range('E28').Value = '  -4.61%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '3.72'
$ws.Range('E29').Value = '  -0.48%  '
$ws.Range('B30').Value = 'Monero'
$ws.Range('C30').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '167.29'
$ws.Range('E30').Value = '  -1.53%  '
$ws.Range('B31').Value = 'Toncoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.07'
$ws.Range('E31').Value = '  -5.71%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '20.48'
$ws.Range('E32').Value = '  -0.74%  '
$ws.Range('B33').Value = 'InternetComputer(DFINITY)'
$ws.Range('C33').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.94'
$ws.Range('E33').Value = '  +8.50%  '
$ws.Range('B34').Value = 'Kaspa'
$ws.Range('C34').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.121'
$ws.Range('E34').Value = '  +0.83%  '
$ws.Range('B35').Value = 'Hedera'
$ws.Range('C35').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0789'
$ws.Range('E35').Value = '  +8.67%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.122'
$ws.Range('E36').Value = '  +0.09%  '
$ws.Range('B37').Value = 'InjectiveProtocol'
$ws.Range('C37').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '26.20'
$ws.Range('E37').Value = '  +1.26%  '
$ws.Range('B38').Value = 'Filecoin'
$ws.Range('C38').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.63'
$ws.Range('E38').Value = '  +0.25%  '
$ws.Range('B39').Value = 'RenderToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '4.14'
$ws.Range('E39').Value = '  +3.67%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0308'
$ws.Range('E40').Value = '  +3.11%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.25'
$ws.Range('E41').Value = '  +0.92%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '12.30'
$ws.Range('E42').Value = '  +2.80%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.69'
$ws.Range('E43').Value = '  -1.30%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.06'
$ws.Range('E44').Value = '  +1.98%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '62.17'
$ws.Range('E45').Value = '  -2.91%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.198'
$ws.Range('E46').Value = '  -1.92%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '8.60'
$ws.Range('E47').Value = '  -0.24%  '
$ws.Range('E48').Value = '  -0.99%  '
$ws.Range('E49').Value = '  -0.29%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.17'
$ws.Range('E50').Value = '  +2.10%  '
$ws.Range('E51').Value = '  +4.51%  '
